# Rewrite the "m:self.someCustomService()" field so that its field code is
# expressed as literal braces text ("{m:self.someCustomService()}") instead
# of a real Word field (fldChar begin/end + instrText runs).
#
# Matches TokenIteratorFieldRewriterSplit's output: every token of the
# former field code becomes its own run (<w:t>), and the "self" token keeps
# the original orange accent color run properties.

$d = $word.ActiveDocument

# Locate the field holding the custom-service instruction text.
$field = $null
foreach ($f in $d.Fields) {
    if ($f.Code.Text -match "someCustomService") {
        $field = $f
        break
    }
}

# Find the paragraph that hosts the field (Field.Code.Paragraphs is not
# reliable in this host, so locate it by scanning Document.Paragraphs for
# the one whose range contains the field's code range).
$fieldPara = $null
foreach ($para in $d.Paragraphs) {
    if ($field.Code.Start -ge $para.Range.Start -and $field.Code.Start -lt $para.Range.End) {
        $fieldPara = $para
        break
    }
}
$insertPos = $fieldPara.Range.Start

# Remove the field (fldChar begin/end + instrText runs) - leaves the
# paragraph empty so we can replace it with plain literal-text runs.
$field.Delete()

$insertAt = $d.Range($insertPos, $insertPos)

# Rebuild the field code as plain text runs, one run per former token, so
# that "self" keeps its own run (and its original orange accent color).
$runsXml = '<w:r><w:t>{</w:t></w:r>' +
           '<w:r><w:t>m</w:t></w:r>' +
           '<w:r><w:t>:</w:t></w:r>' +
           '<w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>self</w:t></w:r>' +
           '<w:r><w:t>.</w:t></w:r>' +
           '<w:r><w:t>someCustomService()</w:t></w:r>' +
           '<w:r><w:t xml:space="preserve">}</w:t></w:r>'

$package = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $runsXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

$insertAt.InsertXML($package) | Out-Null
